# Scheduled-runner market data refresh for Aegis_Profits sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the specific Leve rows whose underlying market prices moved.
$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 566.1212
$ws.Range("I28").Value = 210.1
$ws.Range("J28").Value = 1113.8462
$ws.Range("K28").Value = 210.1
$ws.Range("L28").Value = 1113.8462
$ws.Range("M28").Value = 274.9
$ws.Range("N28").Value = -2083.8462

# ALC row 108
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 26624.5
$ws.Range("J108").Value = 26624.5
$ws.Range("L108").Value = 26624.5
$ws.Range("N108").Value = -34304.5

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1117.2
$ws.Range("I112").Value = 513.3333
$ws.Range("J112").Value = 1199.5454
$ws.Range("K112").Value = 1539.9999
$ws.Range("L112").Value = 3598.6362
$ws.Range("M112").Value = -431.9999
$ws.Range("N112").Value = -5814.6362

# ALC row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3046.5881
$ws.Range("I125").Value = 2602.4614
$ws.Range("J125").Value = 4490
$ws.Range("K125").Value = 23422.1526
$ws.Range("L125").Value = 40410
$ws.Range("M125").Value = -20962.1526
$ws.Range("N125").Value = -45330

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1166880.1
$ws.Range("I129").Value = 50500
$ws.Range("J129").Value = 1390156.2
$ws.Range("K129").Value = 151500
$ws.Range("L129").Value = 4170468.6
$ws.Range("M129").Value = -146500
$ws.Range("N129").Value = -4180468.6

# ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 38836.4
$ws.Range("J133").Value = 38836.4
$ws.Range("L133").Value = 38836.4
$ws.Range("N133").Value = -48956.4

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1430.907
$ws.Range("I137").Value = 1001.2778
$ws.Range("K137").Value = 3003.8334
$ws.Range("M137").Value = -453.8334

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 33513.562
$ws.Range("I45").Value = 51908.15
$ws.Range("J45").Value = 2855.9167
$ws.Range("K45").Value = 51908.15
$ws.Range("L45").Value = 2855.9167
$ws.Range("M45").Value = -51531.15
$ws.Range("N45").Value = -3609.9167

# ARM row 98
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H98").Value = 15999
$ws.Range("J98").Value = 15999
$ws.Range("L98").Value = 15999
$ws.Range("N98").Value = -21989

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2400.4443
$ws.Range("I122").Value = 2400.4443
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7201.3329
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4751.3329
$ws.Range("N122").ClearContents()

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1892.2
$ws.Range("I16").Value = 2122.125
$ws.Range("J16").Value = 972.5
$ws.Range("K16").Value = 2122.125
$ws.Range("L16").Value = 972.5
$ws.Range("M16").Value = -1835.125
$ws.Range("N16").Value = -1546.5

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 19470
$ws.Range("I99").Value = 3264
$ws.Range("J99").Value = 35676
$ws.Range("K99").Value = 3264
$ws.Range("L99").Value = 35676
$ws.Range("M99").Value = -1766
$ws.Range("N99").Value = -38672

# CRP row 104
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 31193.75
$ws.Range("J104").Value = 31193.75
$ws.Range("L104").Value = 31193.75
$ws.Range("N104").Value = -36435.75

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2007.5238
$ws.Range("I105").Value = 1990.8889
$ws.Range("J105").Value = 2107.3333
$ws.Range("K105").Value = 1990.8889
$ws.Range("L105").Value = 2107.3333
$ws.Range("M105").Value = -243.8888999999999
$ws.Range("N105").Value = -5601.3333

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1892.2
$ws.Range("I113").Value = 2122.125
$ws.Range("J113").Value = 972.5
$ws.Range("K113").Value = 2122.125
$ws.Range("L113").Value = 972.5
$ws.Range("M113").Value = 47.875
$ws.Range("N113").Value = -5312.5

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 19470
$ws.Range("I126").Value = 3264
$ws.Range("J126").Value = 35676
$ws.Range("K126").Value = 9792
$ws.Range("L126").Value = 107028
$ws.Range("M126").Value = -7322
$ws.Range("N126").Value = -111968

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 375.55554
$ws.Range("J34").Value = 702.5
$ws.Range("L34").Value = 2107.5
$ws.Range("N34").Value = -2275.5

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1707.7273
$ws.Range("I132").Value = 796.6667
$ws.Range("J132").Value = 2801
$ws.Range("K132").Value = 7170.0003
$ws.Range("L132").Value = 25209
$ws.Range("M132").Value = -4640.0003
$ws.Range("N132").Value = -30269

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2722.1667
$ws.Range("I7").Value = 1658.25
$ws.Range("J7").Value = 4850
$ws.Range("K7").Value = 1658.25
$ws.Range("L7").Value = 4850
$ws.Range("M7").Value = -1546.25
$ws.Range("N7").Value = -5074

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 64356.312
$ws.Range("I40").Value = 167633.67
$ws.Range("J40").Value = 2389.9
$ws.Range("K40").Value = 167633.67
$ws.Range("L40").Value = 2389.9
$ws.Range("M40").Value = -167497.67
$ws.Range("N40").Value = -2661.9

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1125605.8
$ws.Range("I46").Value = 395
$ws.Range("J46").Value = 1447094.6
$ws.Range("K46").Value = 395
$ws.Range("L46").Value = 1447094.6
$ws.Range("M46").Value = -207
$ws.Range("N46").Value = -1447470.6

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 364.77777
$ws.Range("I55").Value = 176.16667
$ws.Range("J55").Value = 553.3889
$ws.Range("K55").Value = 176.16667
$ws.Range("L55").Value = 553.3889
$ws.Range("M55").Value = -3.166670000000011
$ws.Range("N55").Value = -899.3889

# LTW row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 17000
$ws.Range("J106").Value = 17000
$ws.Range("L106").Value = 17000
$ws.Range("N106").Value = -19524

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2722.1667
$ws.Range("I126").Value = 1658.25
$ws.Range("J126").Value = 4850
$ws.Range("K126").Value = 4974.75
$ws.Range("L126").Value = 14550
$ws.Range("M126").Value = -2504.75
$ws.Range("N126").Value = -19490

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1403.5
$ws.Range("I136").Value = 1354.28
$ws.Range("J136").Value = 1579.2858
$ws.Range("K136").Value = 4062.84
$ws.Range("L136").Value = 4737.857400000001
$ws.Range("M136").Value = -1512.84
$ws.Range("N136").Value = -9837.857400000001

# WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5540
$ws.Range("I14").Value = 7700
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 7700
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -7532
$ws.Range("N14").Value = -5336

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1923.1538
$ws.Range("J122").Value = 2714.5715
$ws.Range("L122").Value = 8143.7145
$ws.Range("N122").Value = -13043.7145

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1190.6364
$ws.Range("I126").Value = 1357.9333
$ws.Range("J126").Value = 832.1429000000001
$ws.Range("K126").Value = 4073.7999
$ws.Range("L126").Value = 2496.4287
$ws.Range("M126").Value = -1603.7999
$ws.Range("N126").Value = -7436.4287
